# Jupyter book compile files
# Reorganize the project-topics log on the "ch 1" sheet: the "core_principles"
# row is retired and replaced with a new "authenticity" topic, and the
# "assessment evaluation" section header moves up to row 6 (displacing
# "constructivism"/"connectivism" down one row each).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch 1")

# Row 4 "teaching approach" becomes bold (matches the other section-ish labels)
$ws.Range("B4").Font.Bold = $true

# Row 6 becomes the new "assessment evaluation" header: bold row height, wrap
# text styling (like the row 9 header), and its word count jumps to 872.
$ws.Rows.Item(6).RowHeight = 17
$ws.Range("B6").Value = "assessment evaluation"
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = 872

# Row 7 now holds "constructivism" with its old word count.
$ws.Range("B7").Value = "constructivism"
$ws.Range("C7").Value = 467

# Row 8 now holds "connectivism" (no longer bold - the old "core_principles"
# styling is dropped).
$ws.Range("B8").Value = "connectivism"
$ws.Range("B8").Font.Bold = $false
$ws.Range("C8").Value = 426

# Row 9 keeps its header styling but becomes the new "authenticity" topic.
$ws.Range("B9").Value = "authenticity"
$ws.Range("C9").Value = 413

# Row 10 "open" becomes bold.
$ws.Range("B10").Font.Bold = $true

# Row 14 gains the usual Started?/Finished? "y" markers.
$ws.Range("E14").Value = "y"
$ws.Range("F14").Value = "y"

# Move the active selection to B10 (where the editing ended up).
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
